$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update counter
$ws.Range("F1").Value = 17

# Row 4
$ws.Range("A4").Value = 'ññ'
$ws.Range("B4").Value = '76k1@J0`q1}Z0/Z9&m2@g1/t7>{132k1@130J0`146q1}127Z0/114Z9&153m2@127g1/171t7>{327F3&327M3_{'
$ws.Range("C4").Value = '11m1}L6|R0[w0*O2?)60m1}65L6|132R0[68w0*57O2?)257y9.)'

# Row 5
$ws.Range("A5").Value = 'ñ'
$ws.Range("B5").Value = '62d9$M5<R2{V4(g5!&116d9$110M5<128R2{118V4(125g5!&301B8?&'
$ws.Range("C5").Value = '60d7%O5/T8}L9*L8>:115d7%111O5/162T8}108L9*106L8>:314f0.:'

# Row 6
$ws.Range("A6").Value = 'l'
$ws.Range("B6").Value = '46m1+J3)b8!J1,w5,-101m1+101J3)158b8!97J1,108w5,-185p3>-'
$ws.Range("C6").Value = '95B2*H2^P7.R1&P6|?146B2*151H2^177P7.149R1&156P6|?146R6=?'

# Row 7
$ws.Range("A7").Value = 'camilo'
$ws.Range("B7").Value = '66S8_l5@F8/X0-h1#m4$v7|A5=M0{W4|V8(W0?j5>q5+A4-L7~K8&t7;a4)Q8~[117S8_121l5@149F8/119X0-192h1#132m4$118v7|160A5=146M0{120W4|112V8(177W0?117j5>128q5+167A4-118L7~110K8&174t7;116a4)190Q8~[136S5~134B4^146P6.142o3>145e4,148l2|['
$ws.Range("C7").Value = '87q7*D4&K1!e6$m0?M3&c3_f5!l4&X4>x9*A6_y9:w8;z0]q0|M4.h1$o1:p9>@138q7*144D4&201K1!141e6$132m0?152M3&144c3_213f5!177l4&136X4>120x9*208A6_140y9:150w8;166z0]138q0|178M4.163h1$142o1:150p9>@138r6-136A9~148Z1!144y5?147O3[150L7?@'

# Row 8
$ws.Range("A8").Value = 'camilo'
$ws.Range("B8").Value = '63H0&F6{W4@i5]G6:m3_m8!J2[A1&N0_W4]l6+I8_n2#A0|s4)o1%S2?J1-K2#`118H0&118F6{167W4@117i5]186G6:130m3_116m8!96J2[128A1&118N0_157W4]148l6+113I8_158n2#141A0|111s4)109o1%167S2?118J1-125K2#`176h6{174C5!186A7^182U2_185N0.188h7>`'
$ws.Range("C8").Value = '91b5:Z0;o5?O7|n9?l7`l9}T8)Y6|T3]P6>Q6#Q6~e0(B3%W2,E7|A5}X3}m1;<146b5:145Z0;164o5?147O7|216n9?158l7`147l9}153T8)178Y6|143T3]131P6>181Q6#140Q6~150e0(165B3%145W2,215E7|164A5}141X3}128m1;<175I8}173C8>185W4(181Z1;184J6|187I2%<'

# Row 9
$ws.Range("A9").Value = 'camilo'
$ws.Range("B9").Value = '34D6&Z2_F0?d2{V1?I2#Y8|b9?v6#A6#a2]Q2|c7@s3(h7-w7?b5+r4!N0[F0}/91D6&84Z2_110F0?86d2{74V1?99I2#91Y8|128b9?134v6#86A6#127a2]141Q2|85c7@71s3(144h7-88w7?70b5+133r4!83N0[157F0}/191L4(189A9^201d4]197k3%200n6$203c1{/'
$ws.Range("C9").Value = '30e5`d8_s6|e1[v7^o5;b3`E9.F0~G8<k1!m8~w3*d4-y8*R0=u1(D3,F8]M2:)83e5`83d8_96s6|79e1[121v7^131o5;84b3`70E9.107F0~86G8<125k1!116m8~79w3*156d4-141y8*78R0=88u1(102D3,78F8]91M2:)154B1[152e6(164M8_160V1~163o0:166H0=)'
$ws.Range("B9:C9").Style = "Normal"

# Row 10
$ws.Range("A10").Value = 'a'
$ws.Range("B10").Value = '49i0:g3(K0:m0)H5!=98i0:102g3(155K0:106m0)82H5!=112j9!='
$ws.Range("C10").Value = '24W0@A8+m7:X9#A6~/80W0@74A8+144m7:75X9#67A6~/179x3+/'
$ws.Range("B10:C10").Style = "Normal"

# Row 11
$ws.Range("A11").Value = 'ccc'
$ws.Range("B11").Value = '13d0&L9*A9,g4!S3_w1=D0$q4&<62d0&66L9*113A9,67g4!74S3_102w1=62D0$76q4&<114d6=114Y1?<'
$ws.Range("C11").Value = '27D9[D2)P0}y7*M4.`84D9[82D2)126P0}76y7*89M4.`196c1>`'
$ws.Range("B11:C11").Style = "Normal"

# Row 12
$ws.Range("A12").Value = 'dd'
$ws.Range("B12").Value = '37p8(N1>z4.b7+c3~i5`z1+s8@^86p8(94N1>102z4.93b7+98c3~139i5`94z1+97s8@^134A8=134f9<^'
$ws.Range("C12").Value = '12a1.Q5[m4`c9,Q3`J9(Q0)A3~|65a1.61Q5[123m4`60c9,107Q3`125J9(66Q0)106A3~|148o0_148q6^|'
$ws.Range("B12:C12").Style = "Normal"

# Row 16
$ws.Range("A16").Value = 'camilo'
$ws.Range("B16").Value = '26A8~M3{V1>W9%f4=r9`w8+O0]r1_c8|L7~h2%z2~M4?b6+m1+H5-w3^f6@R4()80A8~76M3{108V1>76W9%63f4=113r9`77w8+84O0]133r1_82c8|87L7~140h2%74z2~73M4?138b6+78m1+66H5-111w3^83f6@117R4()161R2%159W3:171k8=167r0/170p4(173U9[)'
$ws.Range("C16").Value = '46K2#R7*P2-Z5^v4$M2.O6|V3?r4~H1{A2-h4.y8&y3~p5+I7`J0<J5^l3&Z7@!99K2#95R7*165P2-102Z5^140v4$125M2.102O6|86V3?146r4~103H1{172A2-154h4.97y8&91y3~112p5+99I7`142J0<144J5^103l3&171Z7@!150w8^148O8(160d9~156l3-159B5`162b9}!'
$ws.Range("B16:C16").Style = "Normal"
